$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (pushing the existing "Melody toffee" row,
# and everything below it, down by one row).
$ws.Rows("48:48").Insert()

# Populate the newly inserted row with the new "Snacks" item.
$ws.Range("A48").Value = "Snacks"
$ws.Range("B48").Value = "Haldiram's All in one 200 gm"
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 57
$ws.Range("E48").Value = "Haldiram's All in one 200 gm"
$ws.Range("F48").Value = "Fast Food"

# The sheet uses an (inactive) AutoFilter backed by the hidden defined name
# _FilterDatabase. Excel normally grows this range automatically when rows
# are inserted inside it; extend it explicitly from F58 to F59 to mirror
# that behaviour.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$59"
    }
}

# Reflect the author's on-screen scroll position / selection at save time.
[void]$ws.Range("G48").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
